# Refresh cached market-price-derived figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets, per the scheduled price-data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Cells.Item(17, 8).Value = 2775.577
$ws.Cells.Item(17, 10).Value = 2775.577
$ws.Cells.Item(17, 12).Value = 8326.731
$ws.Cells.Item(17, 14).Value = -8662.731

# row 40
$ws.Cells.Item(40, 8).Value = 1722
$ws.Cells.Item(40, 9).Value = 1606.625
$ws.Cells.Item(40, 10).Value = 1985.7142
$ws.Cells.Item(40, 11).Value = 1606.625
$ws.Cells.Item(40, 12).Value = 1985.7142
$ws.Cells.Item(40, 13).Value = -1431.625
$ws.Cells.Item(40, 14).Value = -2335.7142

# row 64
$ws.Cells.Item(64, 8).Value = 6538747
$ws.Cells.Item(64, 9).Value = 11113770
$ws.Cells.Item(64, 10).Value = 2999.1428
$ws.Cells.Item(64, 11).Value = 11113770
$ws.Cells.Item(64, 12).Value = 2999.1428
$ws.Cells.Item(64, 13).Value = -11113522
$ws.Cells.Item(64, 14).Value = -3495.1428

# row 67
$ws.Cells.Item(67, 8).Value = 6538747
$ws.Cells.Item(67, 9).Value = 11113770
$ws.Cells.Item(67, 10).Value = 2999.1428
$ws.Cells.Item(67, 11).Value = 11113770
$ws.Cells.Item(67, 12).Value = 2999.1428
$ws.Cells.Item(67, 13).Value = -11112912
$ws.Cells.Item(67, 14).Value = -4715.1428

# row 69
$ws.Cells.Item(69, 8).Value = 4223.636
$ws.Cells.Item(69, 9).Value = 5125
$ws.Cells.Item(69, 10).Value = 3708.5715
$ws.Cells.Item(69, 11).Value = 15375
$ws.Cells.Item(69, 12).Value = 11125.7145
$ws.Cells.Item(69, 13).Value = -14501
$ws.Cells.Item(69, 14).Value = -12873.7145

# row 72
$ws.Cells.Item(72, 8).Value = 4223.636
$ws.Cells.Item(72, 9).Value = 5125
$ws.Cells.Item(72, 10).Value = 3708.5715
$ws.Cells.Item(72, 11).Value = 46125
$ws.Cells.Item(72, 12).Value = 33377.1435
$ws.Cells.Item(72, 13).Value = -41757
$ws.Cells.Item(72, 14).Value = -42113.1435

# row 97
$ws.Cells.Item(97, 8).Value = 100
$ws.Cells.Item(97, 9).Value = 100
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 300
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 196
$ws.Cells.Item(97, 14).Value = $null

$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Cells.Item(61, 8).Value = 22729692
$ws.Cells.Item(61, 9).Value = 33336238
$ws.Cells.Item(61, 10).Value = 1378.4286
$ws.Cells.Item(61, 11).Value = 33336238
$ws.Cells.Item(61, 12).Value = 1378.4286
$ws.Cells.Item(61, 13).Value = -33336026
$ws.Cells.Item(61, 14).Value = -1802.4286

# row 101
$ws.Cells.Item(101, 8).Value = 51882
$ws.Cells.Item(101, 10).Value = 51882
$ws.Cells.Item(101, 12).Value = 51882
$ws.Cells.Item(101, 14).Value = -58372

# row 133
$ws.Cells.Item(133, 8).Value = 47810.44
$ws.Cells.Item(133, 10).Value = 47810.44
$ws.Cells.Item(133, 12).Value = 47810.44
$ws.Cells.Item(133, 14).Value = -52870.44

# row 136
$ws.Cells.Item(136, 8).Value = 22729692
$ws.Cells.Item(136, 9).Value = 33336238
$ws.Cells.Item(136, 10).Value = 1378.4286
$ws.Cells.Item(136, 11).Value = 100008714
$ws.Cells.Item(136, 12).Value = 4135.2858
$ws.Cells.Item(136, 13).Value = -100006164
$ws.Cells.Item(136, 14).Value = -9235.2858

$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Cells.Item(99, 8).Value = 1014.5455
$ws.Cells.Item(99, 9).Value = 1002
$ws.Cells.Item(99, 11).Value = 1002
$ws.Cells.Item(99, 13).Value = 496

# row 102
$ws.Cells.Item(102, 8).Value = 50049.6

# row 105
$ws.Cells.Item(105, 8).Value = 4300.25
$ws.Cells.Item(105, 9).Value = 3334.4614
$ws.Cells.Item(105, 10).Value = 4846.1304
$ws.Cells.Item(105, 11).Value = 3334.4614
$ws.Cells.Item(105, 12).Value = 4846.1304
$ws.Cells.Item(105, 13).Value = -1587.4614
$ws.Cells.Item(105, 14).Value = -8340.1304

$ws = $wb.Worksheets.Item("CRP")
# row 105
$ws.Cells.Item(105, 8).Value = 948.2
$ws.Cells.Item(105, 9).Value = 933.3333
$ws.Cells.Item(105, 10).Value = 970.5
$ws.Cells.Item(105, 11).Value = 933.3333
$ws.Cells.Item(105, 12).Value = 970.5
$ws.Cells.Item(105, 13).Value = 813.6667
$ws.Cells.Item(105, 14).Value = -4464.5

# row 106
$ws.Cells.Item(106, 8).Value = 50004.332
$ws.Cells.Item(106, 10).Value = 50004.332
$ws.Cells.Item(106, 12).Value = 50004.332
$ws.Cells.Item(106, 14).Value = -52528.332

$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Cells.Item(80, 8).Value = 33338124
$ws.Cells.Item(80, 9).Value = 83335580
$ws.Cells.Item(80, 10).Value = 6483.3335
$ws.Cells.Item(80, 11).Value = 83335580
$ws.Cells.Item(80, 12).Value = 6483.3335
$ws.Cells.Item(80, 13).Value = -83334582
$ws.Cells.Item(80, 14).Value = -8479.333500000001

# row 83
$ws.Cells.Item(83, 8).Value = 33338124
$ws.Cells.Item(83, 9).Value = 83335580
$ws.Cells.Item(83, 10).Value = 6483.3335
$ws.Cells.Item(83, 11).Value = 416677900
$ws.Cells.Item(83, 12).Value = 32416.6675
$ws.Cells.Item(83, 13).Value = -416672908
$ws.Cells.Item(83, 14).Value = -42400.6675

# row 126
$ws.Cells.Item(126, 8).Value = 3913.4443
$ws.Cells.Item(126, 9).Value = 2649.2778
$ws.Cells.Item(126, 10).Value = 5177.6113
$ws.Cells.Item(126, 11).Value = 7947.8334
$ws.Cells.Item(126, 12).Value = 15532.8339
$ws.Cells.Item(126, 13).Value = -5477.8334
$ws.Cells.Item(126, 14).Value = -20472.8339

# row 134
$ws.Cells.Item(134, 8).Value = 26395.6
$ws.Cells.Item(134, 10).Value = 26395.6
$ws.Cells.Item(134, 12).Value = 79186.79999999999
$ws.Cells.Item(134, 14).Value = -84256.79999999999

# row 138
$ws.Cells.Item(138, 8).Value = 58949.25
$ws.Cells.Item(138, 10).Value = 58949.25
$ws.Cells.Item(138, 12).Value = 58949.25
$ws.Cells.Item(138, 14).Value = -69229.25

$ws = $wb.Worksheets.Item("LTW")
# row 81
$ws.Cells.Item(81, 8).Value = 37190.445
$ws.Cells.Item(81, 10).Value = 37190.445
$ws.Cells.Item(81, 12).Value = 37190.445
$ws.Cells.Item(81, 14).Value = -39186.445

# row 84
$ws.Cells.Item(84, 8).Value = 37190.445
$ws.Cells.Item(84, 10).Value = 37190.445
$ws.Cells.Item(84, 12).Value = 111571.335
$ws.Cells.Item(84, 14).Value = -121555.335

# row 100
$ws.Cells.Item(100, 8).Value = 2500.2727
$ws.Cells.Item(100, 9).Value = 2417.1667
$ws.Cells.Item(100, 10).Value = 2600
$ws.Cells.Item(100, 11).Value = 2417.1667
$ws.Cells.Item(100, 12).Value = 2600
$ws.Cells.Item(100, 13).Value = -1876.1667
$ws.Cells.Item(100, 14).Value = -3682

$ws = $wb.Worksheets.Item("WVR")
# row 76
$ws.Cells.Item(76, 8).Value = 15333.333
$ws.Cells.Item(76, 9).Value = 10000
$ws.Cells.Item(76, 10).Value = 18000
$ws.Cells.Item(76, 11).Value = 10000
$ws.Cells.Item(76, 12).Value = 18000
$ws.Cells.Item(76, 13).Value = -9685
$ws.Cells.Item(76, 14).Value = -18630

# row 79
$ws.Cells.Item(79, 8).Value = 15333.333
$ws.Cells.Item(79, 9).Value = 10000
$ws.Cells.Item(79, 10).Value = 18000
$ws.Cells.Item(79, 11).Value = 10000
$ws.Cells.Item(79, 12).Value = 18000
$ws.Cells.Item(79, 13).Value = -8908
$ws.Cells.Item(79, 14).Value = -20184

# row 80
$ws.Cells.Item(80, 8).Value = 38148.832
$ws.Cells.Item(80, 10).Value = 37778.6
$ws.Cells.Item(80, 12).Value = 37778.6
$ws.Cells.Item(80, 14).Value = -39774.6

# row 83
$ws.Cells.Item(83, 8).Value = 38148.832
$ws.Cells.Item(83, 10).Value = 37778.6
$ws.Cells.Item(83, 12).Value = 113335.8
$ws.Cells.Item(83, 14).Value = -123319.8

# row 107
$ws.Cells.Item(107, 8).Value = 1129.4615
$ws.Cells.Item(107, 9).Value = 1571.5
$ws.Cells.Item(107, 10).Value = 422.2
$ws.Cells.Item(107, 11).Value = 4714.5
$ws.Cells.Item(107, 12).Value = 1266.6
$ws.Cells.Item(107, 13).Value = -2794.5
$ws.Cells.Item(107, 14).Value = -5106.6

# row 138
$ws.Cells.Item(138, 8).Value = 41966.668
$ws.Cells.Item(138, 10).Value = 41966.668
$ws.Cells.Item(138, 12).Value = 41966.668
$ws.Cells.Item(138, 14).Value = -52246.668
